$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("TXT_ENHENCE_START", "开始强化"),
    @("TXT_COST", "花费"),
    @("TXT_ENHENCE_CONSUME", "消耗装备"),
    @("TXT_ENHENCE_TITLE", "装备强化"),
    @("TXT_ENHENCE_OK", "确定"),
    @("TXT_ENHENCE_TARGET", "选择强化装备"),
    @("TXT_ENHENCE_CHOOSE", "选择消耗装备")
)

$startRow = 163
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

[void]$ws.Range("B165").Select()
